$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.530.03"
$ws.Range("E2").Value = "  -0.76%  "

$ws.Range("D3").Value = "1.837.51"
$ws.Range("E3").Value = "  -0.77%  "

$ws.Range("D4").Value = "'1.006"
$ws.Range("E4").Value = "  -2.90%  "

$ws.Range("D5").Value = "'316.10"
$ws.Range("E5").Value = "  -2.07%  "

$ws.Range("D6").Value = "'1.005"
$ws.Range("E6").Value = "  -2.77%  "

$ws.Range("D7").Value = "'0.4303"
$ws.Range("E7").Value = "  -2.12%  "

$ws.Range("D8").Value = "'0.3724"
$ws.Range("E8").Value = "  -1.87%  "

$ws.Range("E9").Value = "  -1.60%  "

$ws.Range("D10").Value = "'0.8686"
$ws.Range("E10").Value = "  -1.67%  "

$ws.Range("D11").Value = "'21.25"
$ws.Range("E11").Value = "  -1.40%  "

$ws.Range("D12").Value = "1.846.93"
$ws.Range("E12").Value = "  -0.37%  "

$ws.Range("D13").Value = "'6.706"
$ws.Range("E13").Value = "  +0.04%  "

$ws.Range("D14").Value = "'5.370"
$ws.Range("E14").Value = "  -2.36%  "

$ws.Range("D15").Value = "'0.07082"
$ws.Range("E15").Value = "  -1.26%  "

$ws.Range("D16").Value = "'88.55"
$ws.Range("E16").Value = "  +4.18%  "

$ws.Range("E17").Value = "  -2.91%  "

$ws.Range("D18").Value = "'0.000008945"
$ws.Range("E18").Value = "  -1.41%  "

$ws.Range("D19").Value = "'1.005"
$ws.Range("E19").Value = "  -2.73%  "

$ws.Range("E20").Value = "  -1.12%  "

$ws.Range("D21").Value = "27.539.20"
$ws.Range("E21").Value = "  -0.78%  "

$ws.Range("D22").Value = "'5.170"
$ws.Range("E22").Value = "  -2.28%  "

$ws.Range("D24").Value = "2.071.21"
$ws.Range("E24").Value = "  -0.61%  "

$ws.Range("D25").Value = "'2.009"
$ws.Range("E25").Value = "  -2.32%  "

$ws.Range("D26").Value = "'153.91"
$ws.Range("E26").Value = "  -3.13%  "

$ws.Range("D27").Value = "'18.45"
$ws.Range("E27").Value = "  -1.37%  "

$ws.Range("E28").Value = "  +8.43%  "

$ws.Range("D29").Value = "'5.306"
$ws.Range("E29").Value = "  -0.53%  "

$ws.Range("D30").Value = "'117.37"
$ws.Range("E30").Value = "  -0.36%  "

$ws.Range("D31").Value = "'0.08882"
$ws.Range("E31").Value = "  -2.28%  "

$ws.Range("E32").Value = "  +0.13%  "

$ws.Range("D33").Value = "'0.7722"
$ws.Range("E33").Value = "  -0.01%  "

$ws.Range("D34").Value = "'4.502"
$ws.Range("E34").Value = "  -1.21%  "

$ws.Range("D35").Value = "'2.897"
$ws.Range("E35").Value = "  -3.76%  "

$ws.Range("D36").Value = "'1.006"
$ws.Range("E36").Value = "  -2.79%  "

$ws.Range("D37").Value = "'1.126"
$ws.Range("E37").Value = "  -2.15%  "

$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value = "'0.05292"
$ws.Range("E38").Value = "  +0.55%  "

$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.01960"
$ws.Range("E39").Value = "  -0.71%  "

$ws.Range("D40").Value = "'7.161"
$ws.Range("E40").Value = "  +4.15%  "

$ws.Range("D41").Value = "'2.876"
$ws.Range("E41").Value = "  +0.98%  "

$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").Value = "'0.1680"
$ws.Range("E42").Value = "  +0.67%  "

$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").Value = "'0.5100"
$ws.Range("E43").Value = "  -1.53%  "

$ws.Range("D44").Value = "'8.727"
$ws.Range("E44").Value = "  +0.18%  "

$ws.Range("D45").Value = "'10.61"
$ws.Range("E45").Value = "  -1.25%  "

$ws.Range("D46").Value = "'106.54"
$ws.Range("E46").Value = "  -3.35%  "

$ws.Range("D47").Value = "'0.4728"
$ws.Range("E47").Value = "  +0.69%  "

$ws.Range("D48").Value = "'0.06433"
$ws.Range("E48").Value = "  -2.06%  "

$ws.Range("D49").Value = "'1.005"
$ws.Range("E49").Value = "  -2.93%  "

$ws.Range("D50").Value = "'1.677"
$ws.Range("E50").Value = "  -1.44%  "

$ws.Range("D51").Value = "'1.839"
$ws.Range("E51").Value = "  -2.42%  "

# Reset style on quote-prefixed numeric-text cells to drop the quote-prefix flag
$resetCells = @("D4", "D5", "D6", "D7", "D8", "D10", "D11", "D13", "D14", "D15", "D16", "D18", "D19", "D22", "D25", "D26", "D27", "D29", "D30", "D31", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($c in $resetCells) {
    $ws.Range($c).Style = "Normal"
}
